# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each
#   language sheet's Status column.
# - Narrow the (now shorter) status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status values (previously "Ready for handoff").
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the status columns now that the text is shorter.
$newStatusColWidth = 12.576851254417766

$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newStatusColWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newStatusColWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newStatusColWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newStatusColWidth
